$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.059.60"
$ws.Range("E2").Value = "  -0.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.060.55"
$ws.Range("E3").Value = "  -1.36%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.98"
$ws.Range("E5").Value = "  +0.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.672"
$ws.Range("E6").Value = "  +1.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.33"
$ws.Range("E7").Value = "  +6.23%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.28"
$ws.Range("E9").Value = "  -1.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.387"
$ws.Range("E10").Value = "  +2.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0801"
$ws.Range("E11").Value = "  +6.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.109"
$ws.Range("E12").Value = "  +2.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "16.53"
$ws.Range("E13").Value = "  +10.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.360.36"
$ws.Range("E14").Value = "  -1.35%  "

$ws.Range("E15").Value = "  -2.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.63"
$ws.Range("E16").Value = "  +8.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.057.73"
$ws.Range("E17").Value = "  -1.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.026.49"
$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.49"
$ws.Range("E19").Value = "  +13.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "75.79"
$ws.Range("E20").Value = "  +3.89%  "

$ws.Range("E21").Value = "  +8.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.48"
$ws.Range("E22").Value = "  +4.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.51"
$ws.Range("E23").Value = "  -0.95%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  -2.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.30"
$ws.Range("E26").Value = "  +13.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.57"
$ws.Range("E27").Value = "  -1.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.32"
$ws.Range("E28").Value = "  +1.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.30"
$ws.Range("E29").Value = "  -2.39%  "

$ws.Range("E30").Value = "  +2.23%  "

$ws.Range("E31").Value = "  +5.11%  "

$ws.Range("E32").Value = "  +5.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0623"
$ws.Range("E33").Value = "  +0.33%  "

$ws.Range("E34").Value = "  +7.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0883"
$ws.Range("E35").Value = "  -0.84%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.20%  "

$ws.Range("E37").Value = "  +1.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.76"
$ws.Range("E38").Value = "  -3.29%  "

$ws.Range("E39").Value = "  +18.13%  "

$ws.Range("E40").Value = "  +1.75%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.79"
$ws.Range("E41").Value = "  -1.58%  "

$ws.Range("E42").Value = "  -1.01%  "

$ws.Range("E43").Value = "  -0.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.72"
$ws.Range("E44").Value = "  -1.03%  "

$ws.Range("E45").Value = "  +1.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.64"
$ws.Range("E46").Value = "  +15.14%  "

# Rows 47 and 48: RenderToken / FTXToken swap positions with updated data
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.49"
$ws.Range("E47").Value = "  +5.28%  "

$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.93"
$ws.Range("E48").Value = "  -23.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.293.66"
$ws.Range("E49").Value = "  -2.57%  "

$ws.Range("E50").Value = "  -0.61%  "

$ws.Range("E51").Value = "  -0.84%  "
